$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header: "serial number" -> "serial number 2022"
$ws.Range("C4").Value = "serial number 2022"

# Update B10: "Lower" -> "Lower ( middle of wetland areaish)" and clear D10 (the old note cell)
$ws.Range("B10").Value = "Lower ( middle of wetland areaish)"
$ws.Range("D10").ClearContents()

# Add new D column serial numbers for rows 7,8,9
$ws.Range("D7").Value = 72020442
$ws.Range("D8").Value = 72020437
$ws.Range("D9").Value = 78020435

# Update selection to B10
$ws.Range("B10").Select()
